# Apply corrected stock-report figures: for each affected pair/triple of
# rows (grouped by item), the Item Code (B), Rate (E), Qty (F) and Value (G)
# columns - and in two cases Item Name (C) / Price (D) - were fixed to their
# correct values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B142").Value = 48654
$ws.Range("E142").Value = 38.26
$ws.Range("F142").Value = -1
$ws.Range("G142").Value = -32.02
$ws.Range("B143").Value = 63902
$ws.Range("E143").Value = 34.04
$ws.Range("F143").Value = 2
$ws.Range("G143").Value = 64.04000000000001
$ws.Range("B154").Value = 57756
$ws.Range("F154").Value = -100
$ws.Range("G154").Value = -6644
$ws.Range("B155").Value = 53925
$ws.Range("E155").Value = 79.37
$ws.Range("F155").Value = 1
$ws.Range("G155").Value = 66.44
$ws.Range("B156").Value = 64350
$ws.Range("E156").Value = 70.63
$ws.Range("F156").Value = 101
$ws.Range("G156").Value = 6710.44
$ws.Range("B176").Value = 64329
$ws.Range("E176").Value = 128.32
$ws.Range("F176").Value = 6
$ws.Range("G176").Value = 724.14
$ws.Range("B177").Value = 57552
$ws.Range("E177").Value = 136.86
$ws.Range("F177").Value = -5
$ws.Range("G177").Value = -603.45
$ws.Range("B271").Value = 48706
$ws.Range("E271").Value = 39.8
$ws.Range("F271").Value = -144
$ws.Range("G271").Value = -4795.2
$ws.Range("B272").Value = 64973
$ws.Range("E272").Value = 35.4
$ws.Range("F272").Value = 150
$ws.Range("G272").Value = 4995
$ws.Range("B343").Value = 63571
$ws.Range("F343").Value = 29
$ws.Range("G343").Value = 4160.92
$ws.Range("B344").Value = 63531
$ws.Range("F344").Value = 80
$ws.Range("G344").Value = 11478.4
$ws.Range("B347").Value = 55356
$ws.Range("E347").Value = 54.04
$ws.Range("F347").Value = -158
$ws.Range("G347").Value = -7527.12
$ws.Range("B348").Value = 63510
$ws.Range("E348").Value = 50.66
$ws.Range("F348").Value = 167
$ws.Range("G348").Value = 7955.88
$ws.Range("B367").Value = 63563
$ws.Range("E367").Value = 119.04
$ws.Range("F367").Value = 15
$ws.Range("G367").Value = 1679.4
$ws.Range("B368").Value = 61605
$ws.Range("E368").Value = 133.78
$ws.Range("F368").Value = -13
$ws.Range("G368").Value = -1455.48
$ws.Range("B381").Value = 57817
$ws.Range("F381").Value = 3
$ws.Range("G381").Value = 239.43
$ws.Range("B382").Value = 62865
$ws.Range("F382").Value = 151
$ws.Range("G382").Value = 12051.31
$ws.Range("B423").Value = 63102
$ws.Range("C423").Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range("F423").Value = 36
$ws.Range("G423").Value = 2140.92
$ws.Range("B424").Value = 53082
$ws.Range("C424").Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range("F424").Value = 1
$ws.Range("G424").Value = 59.47
$ws.Range("B528").Value = 58047
$ws.Range("D528").Value = 105.54
$ws.Range("E528").Value = 126.1
$ws.Range("F528").Value = 54
$ws.Range("G528").Value = 5699.16
$ws.Range("B529").Value = 47097
$ws.Range("D529").Value = 112.28
$ws.Range("E529").Value = 134.16
$ws.Range("F529").Value = 15
$ws.Range("G529").Value = 1684.2
$ws.Range("B573").Value = 53602
$ws.Range("E573").Value = 15.69
$ws.Range("F573").Value = -231
$ws.Range("G573").Value = -3037.65
$ws.Range("B574").Value = 65068
$ws.Range("E574").Value = 13.97
$ws.Range("F574").Value = 232
$ws.Range("G574").Value = 3050.8
$ws.Range("B578").Value = 45695
$ws.Range("E578").Value = 23.58
$ws.Range("F578").Value = -36
$ws.Range("G578").Value = -710.28
$ws.Range("B579").Value = 64915
$ws.Range("E579").Value = 20.98
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 789.2
$ws.Range("B585").Value = 64927
$ws.Range("E585").Value = 17.26
$ws.Range("F585").Value = 295
$ws.Range("G585").Value = 4784.9
$ws.Range("B586").Value = 45718
$ws.Range("E586").Value = 19.38
$ws.Range("F586").Value = -294
$ws.Range("G586").Value = -4768.68
$ws.Range("B591").Value = 45709
$ws.Range("E591").Value = 15.69
$ws.Range("F591").Value = -300
$ws.Range("G591").Value = -3945
$ws.Range("B592").Value = 64925
$ws.Range("E592").Value = 13.97
$ws.Range("F592").Value = 302
$ws.Range("G592").Value = 3971.3
$ws.Range("B679").Value = 53319
$ws.Range("E679").Value = 310.64
$ws.Range("F679").Value = -6
$ws.Range("G679").Value = -1643.52
$ws.Range("B680").Value = 64810
$ws.Range("E680").Value = 291.22
$ws.Range("F680").Value = 7
$ws.Range("G680").Value = 1917.44
$ws.Range("B701").Value = 64833
$ws.Range("E701").Value = 34.9
$ws.Range("F701").Value = 99
$ws.Range("G701").Value = 3250.17
$ws.Range("B702").Value = 60025
$ws.Range("E702").Value = 37.22
$ws.Range("F702").Value = -98
$ws.Range("G702").Value = -3217.34
$ws.Range("B712").Value = 64830
$ws.Range("E712").Value = 34.9
$ws.Range("F712").Value = 117
$ws.Range("G712").Value = 3841.11
$ws.Range("B713").Value = 60022
$ws.Range("E713").Value = 37.22
$ws.Range("F713").Value = -113
$ws.Range("G713").Value = -3709.79
$ws.Range("B864").Value = 65079
$ws.Range("E864").Value = 43.44
$ws.Range("F864").Value = 21
$ws.Range("G864").Value = 858.27
$ws.Range("B865").Value = 54751
$ws.Range("E865").Value = 46.34
$ws.Range("F865").Value = -19
$ws.Range("G865").Value = -776.53